$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 303, shifting existing rows 303:337 down to 304:338
$ws.Rows.Item(303).Insert()

# Populate the newly inserted row 303 with the new data record
$ws.Cells.Item(303, 1).Value = 6
$ws.Cells.Item(303, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(303, 3).Value = "Metropolitana"
$ws.Cells.Item(303, 4).Value = 44449
$ws.Cells.Item(303, 5).Value = 13
$ws.Cells.Item(303, 6).Value = 100112003
$ws.Cells.Item(303, 7).Value = "Ajo"
$ws.Cells.Item(303, 8).Value = "Chino"
$ws.Cells.Item(303, 9).Value = "Primera"
$ws.Cells.Item(303, 10).Value = 2200
$ws.Cells.Item(303, 11).Value = 14000
$ws.Cells.Item(303, 12).Value = 14500
$ws.Cells.Item(303, 13).Value = 14295
$ws.Cells.Item(303, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(303, 15).Value = "China"
$ws.Cells.Item(303, 16).Value = 1430
$ws.Cells.Item(303, 17).Value = 10
$ws.Cells.Item(303, 18).Value = "Hortaliza"
